$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark MARSTON row (M17) as completed
$ws.Range("M17").Value = 1

# Fill in additional ratings for SUKUP 0022 (row 19)
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 7.75
$ws.Range("F19").Value = 5.5

# Update the active selection to reflect the last edited cell area
$ws.Range("M18").Select()
